# Update Name of Algo
# Apply updated KNN-imputed values in column C for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C4"   = -12.291
    "C7"   = -12.909
    "C16"  = -13.383
    "C28"  = -13.219
    "C29"  = -11.736
    "C32"  = -13.372
    "C40"  = -12.782
    "C52"  = -11.614
    "C57"  = -13.829
    "C66"  = -11.053
    "C100" = -13.322
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
